$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A holds a text-like date string ("01-10-2021"). A plain .Value
# assignment gets auto-parsed by Excel into a date serial + date number
# format, which is not what the source data uses (it stores these as plain
# shared strings with the default/general style, like the rows above it).
# Entering it as a formula that evaluates to the text, then freezing it back
# into a literal via copy / paste-special-values, keeps it as plain text
# without Excel inventing a new date number format / style.
$ws.Range("A47").Formula = "=""01-10-2021"""
$ws.Range("A47").Copy($ws.Range("A47"))
$ws.Range("A47").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B47").Value = 118.03
$ws.Range("C47").Value = 108.87
$ws.Range("D47").Value = 97.11
$ws.Range("E47").Value = 114.26
$ws.Range("F47").Value = 111.31
$ws.Range("G47").Value = 107.8
$ws.Range("H47").Value = 118.79
$ws.Range("I47").Value = 93.32
$ws.Range("J47").Value = 118.26
$ws.Range("K47").Value = 112.31
$ws.Range("L47").Value = 114.55
$ws.Range("M47").Value = 112.37
